$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update kitPartNumber values for rows 7-13 and clear their cell style
# (revert to default "Normal" style, removing the custom style index)
$updates = @{
    7  = 100038
    8  = 100080
    9  = 100081
    10 = 100112
    11 = 100116
    12 = 100118
    13 = 100120
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}

# Update the active selection on the sheet
$ws.Range("C4").Select()
